# Sample Project / Main.xlsx "save" edit:
# Rule cell B11 ("Rule" column, row for the R40/22-23 band) changes from
# the shared text "R40" to the shared text "1". It must remain a plain
# text value (not a number, not a formula) so the cell keeps its existing
# "General" style instead of picking up a new number-format style.
#
# Writing the literal string "1" straight into Range.Value would be
# auto-coerced to a numeric value by Excel, so instead we enter it as a
# text formula and then paste-special just the values back over itself,
# which collapses the formula to a plain (shared-string) text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)
